# Refresh market-price / profit figures on the leve-profit sheets.
# Each block below corresponds to one leve row whose priced columns
# (H..N: currentAveragePrice*, Leve Price*, LeveProfit*) were recomputed
# by the scheduled pricing runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17: H17, J17, L17, N17
$ws.Range("H17").Value = 965.42426
$ws.Range("J17").Value = 965.42426
$ws.Range("L17").Value = 2896.27278
$ws.Range("N17").Value = -3232.27278

# Row 58: H58, I58, K58, M58
$ws.Range("H58").Value = 168
$ws.Range("I58").Value = 168
$ws.Range("K58").Value = 504
$ws.Range("M58").Value = -354

# Row 61: H61, I61, K61, M61
$ws.Range("H61").Value = 304.25
$ws.Range("I61").Value = 304.25
$ws.Range("K61").Value = 912.75
$ws.Range("M61").Value = -740.75

# Row 82: H82, J82, L82, N82
$ws.Range("H82").Value = 7902
$ws.Range("J82").Value = 9942.857
$ws.Range("L82").Value = 29828.571
$ws.Range("N82").Value = -30640.571

# Row 85: H85, J85, L85, N85
$ws.Range("H85").Value = 7902
$ws.Range("J85").Value = 9942.857
$ws.Range("L85").Value = 29828.571
$ws.Range("N85").Value = -32636.571

# Row 97: H97, J97, L97, N97
$ws.Range("H97").Value = 1589.091
$ws.Range("J97").Value = 1648
$ws.Range("L97").Value = 4944
$ws.Range("N97").Value = -5936

# Row 99: H99, I99, K99, M99
$ws.Range("H99").Value = 7129
$ws.Range("I99").Value = 2113.6
$ws.Range("K99").Value = 6340.799999999999
$ws.Range("M99").Value = -4842.799999999999

# Row 101: H101, J101, L101, N101
$ws.Range("H101").Value = 3149.5
$ws.Range("J101").Value = 6047
$ws.Range("L101").Value = 18141
$ws.Range("N101").Value = -21385

# Row 104: H104, I104, K104, M104
$ws.Range("H104").Value = 549.25
$ws.Range("I104").Value = 549.25
$ws.Range("K104").Value = 1647.75
$ws.Range("M104").Value = 99.25

# Row 112: H112, I112, J112, K112, L112, M112, N112
$ws.Range("H112").Value = 4359.6
$ws.Range("I112").Value = 3400
$ws.Range("J112").Value = 4399.5835
$ws.Range("K112").Value = 10200
$ws.Range("L112").Value = 13198.7505
$ws.Range("M112").Value = -9092
$ws.Range("N112").Value = -15414.7505

# Row 115: H115
$ws.Range("H115").Value = 677

# Row 118: H118, I118, K118, M118
$ws.Range("H118").Value = 409
$ws.Range("I118").Value = 409
$ws.Range("K118").Value = 1227
$ws.Range("M118").Value = 430

# Row 127: H127, I127, K127, M127
$ws.Range("H127").Value = 2744.7778
$ws.Range("I127").Value = 1537.5385
$ws.Range("K127").Value = 4612.6155
$ws.Range("M127").Value = 347.3845000000001

# Row 129: H129, I129, J129, K129, L129, M129, N129
$ws.Range("H129").Value = 2294.1177
$ws.Range("I129").Value = 701.1818
$ws.Range("J129").Value = 5214.5
$ws.Range("K129").Value = 2103.5454
$ws.Range("L129").Value = 15643.5
$ws.Range("M129").Value = 2896.4546
$ws.Range("N129").Value = -25643.5

# Row 132: H132, I132, J132, K132, L132, M132, N132
$ws.Range("H132").Value = 4892.4346
$ws.Range("I132").Value = 1958.375
$ws.Range("J132").Value = 11598.857
$ws.Range("K132").Value = 5875.125
$ws.Range("L132").Value = 34796.571
$ws.Range("M132").Value = -3345.125
$ws.Range("N132").Value = -39856.571

# Row 138: H138, I138, J138, K138, L138, M138, N138
$ws.Range("H138").Value = 7477.184
$ws.Range("I138").Value = 2444
$ws.Range("J138").Value = 7908.6
$ws.Range("K138").Value = 7332
$ws.Range("L138").Value = 23725.8
$ws.Range("M138").Value = -2192
$ws.Range("N138").Value = -34005.8

# Row 139: H139, J139, L139, N139
$ws.Range("H139").Value = 113963.164
$ws.Range("J139").Value = 113963.164
$ws.Range("L139").Value = 113963.164
$ws.Range("N139").Value = -124243.164

# Row 140: H140, J140, L140, N140
$ws.Range("H140").Value = 55182.848
$ws.Range("J140").Value = 53889
$ws.Range("L140").Value = 53889
$ws.Range("N140").Value = -64249

$ws = $wb.Worksheets.Item("ARM")
# Row 63: H63, J63, L63, N63
$ws.Range("H63").Value = 3593.125
$ws.Range("J63").Value = 3433.3333
$ws.Range("L63").Value = 3433.3333
$ws.Range("N63").Value = -4805.3333

# Row 66: H66, J66, L66, N66
$ws.Range("H66").Value = 3593.125
$ws.Range("J66").Value = 3433.3333
$ws.Range("L66").Value = 17166.6665
$ws.Range("N66").Value = -24030.6665

# Row 74: H74, I74, J74, K74, L74, M74, N74
$ws.Range("H74").Value = 1352.1428
$ws.Range("I74").Value = 1099.5
$ws.Range("J74").Value = 1903.3636
$ws.Range("K74").Value = 1099.5
$ws.Range("L74").Value = 1903.3636
$ws.Range("M74").Value = -225.5
$ws.Range("N74").Value = -3651.3636

# Row 77: H77, I77, J77, K77, L77, M77, N77
$ws.Range("H77").Value = 1352.1428
$ws.Range("I77").Value = 1099.5
$ws.Range("J77").Value = 1903.3636
$ws.Range("K77").Value = 5497.5
$ws.Range("L77").Value = 9516.818
$ws.Range("M77").Value = -1129.5
$ws.Range("N77").Value = -18252.818

# Row 140: H140, J140, L140, N140
$ws.Range("H140").Value = 115000
$ws.Range("J140").Value = 115000
$ws.Range("L140").Value = 115000
$ws.Range("N140").Value = -125360

$ws = $wb.Worksheets.Item("CRP")
# Row 141: H141, J141, L141, N141
$ws.Range("H141").Value = 111931.7
$ws.Range("J141").Value = 130311.75
$ws.Range("L141").Value = 130311.75
$ws.Range("N141").Value = -140671.75

$ws = $wb.Worksheets.Item("CUL")
# Row 81: H81, J81, L81, N81
$ws.Range("H81").Value = 2767.8635
$ws.Range("J81").Value = 3326.6
$ws.Range("L81").Value = 9979.8
$ws.Range("N81").Value = -12225.8

# Row 84: H84, J84, L84, N84
$ws.Range("H84").Value = 2767.8635
$ws.Range("J84").Value = 3326.6
$ws.Range("L84").Value = 29939.4
$ws.Range("N84").Value = -41171.39999999999

$ws = $wb.Worksheets.Item("GSM")
# Row 2: H2, I2, J2, K2, L2, M2, N2
$ws.Range("H2").Value = 480.53333
$ws.Range("I2").Value = 312.2
$ws.Range("J2").Value = 564.7
$ws.Range("K2").Value = 312.2
$ws.Range("L2").Value = 564.7
$ws.Range("M2").Value = -199.2
$ws.Range("N2").Value = -790.7

# Row 43: H43, I43, J43, K43, L43, M43
$ws.Range("H43").Value = 4997
$ws.Range("I43").Value = 4997
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 4997
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -4846

# Row 46: H46, J46, L46, N46
$ws.Range("H46").Value = 22996.666
$ws.Range("J46").Value = 23497
$ws.Range("L46").Value = 23497
$ws.Range("N46").Value = -23809

# Row 57: H57, I57, J57, K57, L57, M57
$ws.Range("H57").Value = 17500
$ws.Range("I57").Value = 17500
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 17500
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -16680

# Row 113: H113, I113, K113, M113
$ws.Range("H113").Value = 2051
$ws.Range("I113").Value = 2006.25
$ws.Range("K113").Value = 2006.25
$ws.Range("M113").Value = 163.75

$ws = $wb.Worksheets.Item("LTW")
# Row 68: H68, I68, K68, M68
$ws.Range("H68").Value = 1424186.2
$ws.Range("I68").Value = 1896465.5
$ws.Range("K68").Value = 1896465.5
$ws.Range("M68").Value = -1895716.5

# Row 71: H71, I71, K71, M71
$ws.Range("H71").Value = 1424186.2
$ws.Range("I71").Value = 1896465.5
$ws.Range("K71").Value = 9482327.5
$ws.Range("M71").Value = -9478583.5

# Row 82: H82, I82, J82, K82, L82, M82, N82
$ws.Range("H82").Value = 3473411
$ws.Range("I82").Value = 4465542.5
$ws.Range("J82").Value = 949.5
$ws.Range("K82").Value = 4465542.5
$ws.Range("L82").Value = 949.5
$ws.Range("M82").Value = -4465181.5
$ws.Range("N82").Value = -1671.5

# Row 85: H85, I85, J85, K85, L85, M85, N85
$ws.Range("H85").Value = 3473411
$ws.Range("I85").Value = 4465542.5
$ws.Range("J85").Value = 949.5
$ws.Range("K85").Value = 4465542.5
$ws.Range("L85").Value = 949.5
$ws.Range("M85").Value = -4464294.5
$ws.Range("N85").Value = -3445.5

# Row 132: H132, I132, J132, K132, L132, M132, N132
$ws.Range("H132").Value = 3765.4429
$ws.Range("I132").Value = 2815.94
$ws.Range("J132").Value = 6139.2
$ws.Range("K132").Value = 8447.82
$ws.Range("L132").Value = 18417.6
$ws.Range("M132").Value = -5917.82
$ws.Range("N132").Value = -23477.6

$ws = $wb.Worksheets.Item("WVR")
# Row 62: H62, J62, L62, N62
$ws.Range("H62").Value = 7750.125
$ws.Range("J62").Value = 7142.7144
$ws.Range("L62").Value = 7142.7144
$ws.Range("N62").Value = -8390.7144

# Row 65: H65, J65, L65, N65
$ws.Range("H65").Value = 7750.125
$ws.Range("J65").Value = 7142.7144
$ws.Range("L65").Value = 35713.572
$ws.Range("N65").Value = -41953.572

# Row 122: H122, I122, J122, K122, L122, M122, N122
$ws.Range("H122").Value = 4065.9375
$ws.Range("I122").Value = 3950.7036
$ws.Range("J122").Value = 4688.2
$ws.Range("K122").Value = 11852.1108
$ws.Range("L122").Value = 14064.6
$ws.Range("M122").Value = -9402.110799999999
$ws.Range("N122").Value = -18964.6

# Row 132: H132, I132, K132, M132
$ws.Range("H132").Value = 12628668
$ws.Range("I132").Value = 1853701.1
$ws.Range("K132").Value = 5561103.300000001
$ws.Range("M132").Value = -5558573.300000001

# A couple of rows lost their HQ-profit figure (column N) entirely and
# gained (or kept) an NQ-profit figure (column M) instead -- clear the
# now-stale N cells so the row shape matches the refreshed data.
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N43").ClearContents()
$ws.Range("N57").ClearContents()
